$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 7).Value = 13.458797
$ws.Cells.Item(2, 8).Value = 40.376391
$ws.Cells.Item(2, 9).Value = 0.06830096976102129
$ws.Cells.Item(2, 10).Value = 0.06973720484213804
$ws.Cells.Item(2, 13).Value = 38.10639333333333
$ws.Cells.Item(2, 14).Value = 114.31918
$ws.Cells.Item(2, 15).Value = 0.3831479157160237
$ws.Cells.Item(2, 16).Value = 0.4159903984418967
$ws.Cells.Item(2, 17).Value = 512.8662122754865
$ws.Cells.Item(2, 18).Value = 4615.795910479379
$ws.Cells.Item(2, 19).Value = 0.02616937420531847
$ws.Cells.Item(2, 20).Value = 0.02901000762850517
# Row 3
$ws.Cells.Item(3, 7).Value = 13.458797
$ws.Cells.Item(3, 8).Value = 40.376391
$ws.Cells.Item(3, 9).Value = 0.06830096976102129
$ws.Cells.Item(3, 10).Value = 0.06973720484213804
$ws.Cells.Item(3, 15).Value = 0.09199521176963764
$ws.Cells.Item(3, 16).Value = 0.09988081163714851
$ws.Cells.Item(3, 17).Value = 123.141047810748
$ws.Cells.Item(3, 18).Value = 1108.269430296732
$ws.Cells.Item(3, 19).Value = 0.00628336217723677
$ws.Cells.Item(3, 20).Value = 0.006965408620938831
# Row 4
$ws.Cells.Item(4, 7).Value = 13.458797
$ws.Cells.Item(4, 8).Value = 40.376391
$ws.Cells.Item(4, 9).Value = 0.06830096976102129
$ws.Cells.Item(4, 10).Value = 0.06973720484213804
$ws.Cells.Item(4, 13).Value = 15.023598
$ws.Cells.Item(4, 14).Value = 45.070794
$ws.Cells.Item(4, 15).Value = 0.1510575983904562
$ws.Cells.Item(4, 16).Value = 0.1640058785774412
$ws.Cells.Item(4, 17).Value = 202.199555691606
$ws.Cells.Item(4, 18).Value = 1819.796001224454
$ws.Cells.Item(4, 19).Value = 0.01031738045983905
$ws.Cells.Item(4, 20).Value = 0.01143731154966984
# Row 5
$ws.Cells.Item(5, 7).Value = 13.458797
$ws.Cells.Item(5, 8).Value = 40.376391
$ws.Cells.Item(5, 9).Value = 0.06830096976102129
$ws.Cells.Item(5, 10).Value = 0.06973720484213804
$ws.Cells.Item(5, 13).Value = 23.556204
$ws.Cells.Item(5, 14).Value = 47.112408
$ws.Cells.Item(5, 15).Value = 0.236850294013169
$ws.Cells.Item(5, 16).Value = 0.1714350065796238
$ws.Cells.Item(5, 17).Value = 317.038167726588
$ws.Cells.Item(5, 18).Value = 1902.229006359528
$ws.Cells.Item(5, 19).Value = 0.01617710476928246
$ws.Cells.Item(5, 20).Value = 0.01195539817095651
# Row 6
$ws.Cells.Item(6, 7).Value = 13.458797
$ws.Cells.Item(6, 8).Value = 40.376391
$ws.Cells.Item(6, 9).Value = 0.06830096976102129
$ws.Cells.Item(6, 10).Value = 0.06973720484213804
$ws.Cells.Item(6, 13).Value = 13.62041
$ws.Cells.Item(6, 14).Value = 40.86123000000001
$ws.Cells.Item(6, 15).Value = 0.1369489801107134
$ws.Cells.Item(6, 16).Value = 0.1486879047638899
$ws.Cells.Item(6, 17).Value = 183.31433324677
$ws.Cells.Item(6, 18).Value = 1649.82899922093
$ws.Cells.Item(6, 19).Value = 0.00935374814934454
$ws.Cells.Item(6, 20).Value = 0.0103690788720677
# Row 7
$ws.Cells.Item(7, 9).Value = 0.1240039124627887
$ws.Cells.Item(7, 10).Value = 0.1266114708898203
$ws.Cells.Item(7, 13).Value = 38.10639333333333
$ws.Cells.Item(7, 14).Value = 114.31918
$ws.Cells.Item(7, 15).Value = 0.3831479157160237
$ws.Cells.Item(7, 16).Value = 0.4159903984418967
$ws.Cells.Item(7, 17).Value = 931.1349035694932
$ws.Cells.Item(7, 18).Value = 8380.214132125438
$ws.Cells.Item(7, 19).Value = 0.04751184060074975
$ws.Cells.Item(7, 20).Value = 0.05266915622277096
# Row 8
$ws.Cells.Item(8, 9).Value = 0.1240039124627887
$ws.Cells.Item(8, 10).Value = 0.1266114708898203
$ws.Cells.Item(8, 15).Value = 0.09199521176963764
$ws.Cells.Item(8, 16).Value = 0.09988081163714851
$ws.Cells.Item(8, 19).Value = 0.01140776618727785
$ws.Cells.Item(8, 20).Value = 0.01264605647504846
# Row 9
$ws.Cells.Item(9, 9).Value = 0.1240039124627887
$ws.Cells.Item(9, 10).Value = 0.1266114708898203
$ws.Cells.Item(9, 13).Value = 15.023598
$ws.Cells.Item(9, 14).Value = 45.070794
$ws.Cells.Item(9, 15).Value = 0.1510575983904562
$ws.Cells.Item(9, 16).Value = 0.1640058785774412
$ws.Cells.Item(9, 17).Value = 367.103660339328
$ws.Cells.Item(9, 18).Value = 3303.932943053952
$ws.Cells.Item(9, 19).Value = 0.01873173320764922
$ws.Cells.Item(9, 20).Value = 0.02076502552126711
# Row 10
$ws.Cells.Item(10, 9).Value = 0.1240039124627887
$ws.Cells.Item(10, 10).Value = 0.1266114708898203
$ws.Cells.Item(10, 13).Value = 23.556204
$ws.Cells.Item(10, 14).Value = 47.112408
$ws.Cells.Item(10, 15).Value = 0.236850294013169
$ws.Cells.Item(10, 16).Value = 0.1714350065796238
$ws.Cells.Item(10, 17).Value = 575.599048383744
$ws.Cells.Item(10, 18).Value = 3453.594290302464
$ws.Cells.Item(10, 19).Value = 0.02937036312559477
$ws.Cells.Item(10, 20).Value = 0.0217056383450522
# Row 11
$ws.Cells.Item(11, 9).Value = 0.1240039124627887
$ws.Cells.Item(11, 10).Value = 0.1266114708898203
$ws.Cells.Item(11, 13).Value = 13.62041
$ws.Cells.Item(11, 14).Value = 40.86123000000001
$ws.Cells.Item(11, 15).Value = 0.1369489801107134
$ws.Cells.Item(11, 16).Value = 0.1486879047638899
$ws.Cells.Item(11, 17).Value = 332.81657072576
$ws.Cells.Item(11, 18).Value = 2995.34913653184
$ws.Cells.Item(11, 19).Value = 0.01698220934151709
$ws.Cells.Item(11, 20).Value = 0.01882559432568162
# Row 12
$ws.Cells.Item(12, 7).Value = 76.51423666666666
$ws.Cells.Item(12, 8).Value = 229.54271
$ws.Cells.Item(12, 9).Value = 0.3882959647030583
$ws.Cells.Item(12, 10).Value = 0.3964610652618627
$ws.Cells.Item(12, 13).Value = 38.10639333333333
$ws.Cells.Item(12, 14).Value = 114.31918
$ws.Cells.Item(12, 15).Value = 0.3831479157160237
$ws.Cells.Item(12, 16).Value = 0.4159903984418967
$ws.Cells.Item(12, 17).Value = 2915.681598019755
$ws.Cells.Item(12, 18).Value = 26241.1343821778
$ws.Cells.Item(12, 19).Value = 0.1487747895569195
$ws.Cells.Item(12, 20).Value = 0.1649239965049811
# Row 13
$ws.Cells.Item(13, 7).Value = 76.51423666666666
$ws.Cells.Item(13, 8).Value = 229.54271
$ws.Cells.Item(13, 9).Value = 0.3882959647030583
$ws.Cells.Item(13, 10).Value = 0.3964610652618627
$ws.Cells.Item(13, 15).Value = 0.09199521176963764
$ws.Cells.Item(13, 16).Value = 0.09988081163714851
$ws.Cells.Item(13, 17).Value = 700.0657841538799
$ws.Cells.Item(13, 18).Value = 6300.59205738492
$ws.Cells.Item(13, 19).Value = 0.03572136950215359
$ws.Cells.Item(13, 20).Value = 0.03959885298088335
# Row 14
$ws.Cells.Item(14, 7).Value = 76.51423666666666
$ws.Cells.Item(14, 8).Value = 229.54271
$ws.Cells.Item(14, 9).Value = 0.3882959647030583
$ws.Cells.Item(14, 10).Value = 0.3964610652618627
$ws.Cells.Item(14, 13).Value = 15.023598
$ws.Cells.Item(14, 14).Value = 45.070794
$ws.Cells.Item(14, 15).Value = 0.1510575983904562
$ws.Cells.Item(14, 16).Value = 0.1640058785774412
$ws.Cells.Item(14, 17).Value = 1149.51913295686
$ws.Cells.Item(14, 18).Value = 10345.67219661174
$ws.Cells.Item(14, 19).Value = 0.05865505589274934
$ws.Cells.Item(14, 20).Value = 0.06502194533002005
# Row 15
$ws.Cells.Item(15, 7).Value = 76.51423666666666
$ws.Cells.Item(15, 8).Value = 229.54271
$ws.Cells.Item(15, 9).Value = 0.3882959647030583
$ws.Cells.Item(15, 10).Value = 0.3964610652618627
$ws.Cells.Item(15, 13).Value = 23.556204
$ws.Cells.Item(15, 14).Value = 47.112408
$ws.Cells.Item(15, 15).Value = 0.236850294013169
$ws.Cells.Item(15, 16).Value = 0.1714350065796238
$ws.Cells.Item(15, 17).Value = 1802.38496782428
$ws.Cells.Item(15, 18).Value = 10814.30980694568
$ws.Cells.Item(15, 19).Value = 0.09196801340404645
$ws.Cells.Item(15, 20).Value = 0.06796730533173211
# Row 16
$ws.Cells.Item(16, 7).Value = 76.51423666666666
$ws.Cells.Item(16, 8).Value = 229.54271
$ws.Cells.Item(16, 9).Value = 0.3882959647030583
$ws.Cells.Item(16, 10).Value = 0.3964610652618627
$ws.Cells.Item(16, 13).Value = 13.62041
$ws.Cells.Item(16, 14).Value = 40.86123000000001
$ws.Cells.Item(16, 15).Value = 0.1369489801107134
$ws.Cells.Item(16, 16).Value = 0.1486879047638899
$ws.Cells.Item(16, 17).Value = 1042.155274237033
$ws.Cells.Item(16, 18).Value = 9379.397468133302
$ws.Cells.Item(16, 19).Value = 0.0531767363471894
$ws.Cells.Item(16, 20).Value = 0.05894896511424617
# Row 17
$ws.Cells.Item(17, 7).Value = 12.174794
$ws.Cells.Item(17, 8).Value = 24.349588
$ws.Cells.Item(17, 9).Value = 0.06178488588843889
$ws.Cells.Item(17, 10).Value = 0.04205606702633888
$ws.Cells.Item(17, 13).Value = 38.10639333333333
$ws.Cells.Item(17, 14).Value = 114.31918
$ws.Cells.Item(17, 15).Value = 0.3831479157160237
$ws.Cells.Item(17, 16).Value = 0.4159903984418967
$ws.Cells.Item(17, 17).Value = 463.9374889163066
$ws.Cells.Item(17, 18).Value = 2783.62493349784
$ws.Cells.Item(17, 19).Value = 0.02367275025090773
$ws.Cells.Item(17, 20).Value = 0.01749492007918583
# Row 18
$ws.Cells.Item(18, 7).Value = 12.174794
$ws.Cells.Item(18, 8).Value = 24.349588
$ws.Cells.Item(18, 9).Value = 0.06178488588843889
$ws.Cells.Item(18, 10).Value = 0.04205606702633888
$ws.Cells.Item(18, 15).Value = 0.09199521176963764
$ws.Cells.Item(18, 16).Value = 0.09988081163714851
$ws.Cells.Item(18, 17).Value = 111.393082906296
$ws.Cells.Item(18, 18).Value = 668.358497437776
$ws.Cells.Item(18, 19).Value = 0.005683913661469832
$ws.Cells.Item(18, 20).Value = 0.004200594108857047
# Row 19
$ws.Cells.Item(19, 7).Value = 12.174794
$ws.Cells.Item(19, 8).Value = 24.349588
$ws.Cells.Item(19, 9).Value = 0.06178488588843889
$ws.Cells.Item(19, 10).Value = 0.04205606702633888
$ws.Cells.Item(19, 13).Value = 15.023598
$ws.Cells.Item(19, 14).Value = 45.070794
$ws.Cells.Item(19, 15).Value = 0.1510575983904562
$ws.Cells.Item(19, 16).Value = 0.1640058785774412
$ws.Cells.Item(19, 17).Value = 182.909210788812
$ws.Cells.Item(19, 18).Value = 1097.455264732872
$ws.Cells.Item(19, 19).Value = 0.009333076479135964
$ws.Cells.Item(19, 20).Value = 0.006897442222166465
# Row 20
$ws.Cells.Item(20, 7).Value = 12.174794
$ws.Cells.Item(20, 8).Value = 24.349588
$ws.Cells.Item(20, 9).Value = 0.06178488588843889
$ws.Cells.Item(20, 10).Value = 0.04205606702633888
$ws.Cells.Item(20, 13).Value = 23.556204
$ws.Cells.Item(20, 14).Value = 47.112408
$ws.Cells.Item(20, 15).Value = 0.236850294013169
$ws.Cells.Item(20, 16).Value = 0.1714350065796238
$ws.Cells.Item(20, 17).Value = 286.791931121976
$ws.Cells.Item(20, 18).Value = 1147.167724487904
$ws.Cells.Item(20, 19).Value = 0.01463376838824685
$ws.Cells.Item(20, 20).Value = 0.007209882127373508
# Row 21
$ws.Cells.Item(21, 7).Value = 12.174794
$ws.Cells.Item(21, 8).Value = 24.349588
$ws.Cells.Item(21, 9).Value = 0.06178488588843889
$ws.Cells.Item(21, 10).Value = 0.04205606702633888
$ws.Cells.Item(21, 13).Value = 13.62041
$ws.Cells.Item(21, 14).Value = 40.86123000000001
$ws.Cells.Item(21, 15).Value = 0.1369489801107134
$ws.Cells.Item(21, 16).Value = 0.1486879047638899
$ws.Cells.Item(21, 17).Value = 165.82568594554
$ws.Cells.Item(21, 18).Value = 994.9541156732402
$ws.Cells.Item(21, 19).Value = 0.008461377108678513
$ws.Cells.Item(21, 20).Value = 0.006253228488756046
# Row 22
$ws.Cells.Item(22, 7).Value = 70.46836733333333
$ws.Cells.Item(22, 8).Value = 211.405102
$ws.Cells.Item(22, 9).Value = 0.3576142671846927
$ws.Cells.Item(22, 10).Value = 0.36513419197984
$ws.Cells.Item(22, 13).Value = 38.10639333333333
$ws.Cells.Item(22, 14).Value = 114.31918
$ws.Cells.Item(22, 15).Value = 0.3831479157160237
$ws.Cells.Item(22, 16).Value = 0.4159903984418967
$ws.Cells.Item(22, 17).Value = 2685.295323161818
$ws.Cells.Item(22, 18).Value = 24167.65790845636
$ws.Cells.Item(22, 19).Value = 0.1370191611021283
$ws.Cells.Item(22, 20).Value = 0.1518923180064536
# Row 23
$ws.Cells.Item(23, 7).Value = 70.46836733333333
$ws.Cells.Item(23, 8).Value = 211.405102
$ws.Cells.Item(23, 9).Value = 0.3576142671846927
$ws.Cells.Item(23, 10).Value = 0.36513419197984
$ws.Cells.Item(23, 15).Value = 0.09199521176963764
$ws.Cells.Item(23, 16).Value = 0.09988081163714851
$ws.Cells.Item(23, 17).Value = 644.7491994224559
$ws.Cells.Item(23, 18).Value = 5802.742794802104
$ws.Cells.Item(23, 19).Value = 0.03289880024149958
$ws.Cells.Item(23, 20).Value = 0.03646989945142082
# Row 24
$ws.Cells.Item(24, 7).Value = 70.46836733333333
$ws.Cells.Item(24, 8).Value = 211.405102
$ws.Cells.Item(24, 9).Value = 0.3576142671846927
$ws.Cells.Item(24, 10).Value = 0.36513419197984
$ws.Cells.Item(24, 13).Value = 15.023598
$ws.Cells.Item(24, 14).Value = 45.070794
$ws.Cells.Item(24, 15).Value = 0.1510575983904562
$ws.Cells.Item(24, 16).Value = 0.1640058785774412
$ws.Cells.Item(24, 17).Value = 1058.688422532332
$ws.Cells.Item(24, 18).Value = 9528.195802790988
$ws.Cells.Item(24, 19).Value = 0.05402035235108261
$ws.Cells.Item(24, 20).Value = 0.05988415395431776
# Row 25
$ws.Cells.Item(25, 7).Value = 70.46836733333333
$ws.Cells.Item(25, 8).Value = 211.405102
$ws.Cells.Item(25, 9).Value = 0.3576142671846927
$ws.Cells.Item(25, 10).Value = 0.36513419197984
$ws.Cells.Item(25, 13).Value = 23.556204
$ws.Cells.Item(25, 14).Value = 47.112408
$ws.Cells.Item(25, 15).Value = 0.236850294013169
$ws.Cells.Item(25, 16).Value = 0.1714350065796238
$ws.Cells.Item(25, 17).Value = 1659.967236450936
$ws.Cells.Item(25, 18).Value = 9959.803418705616
$ws.Cells.Item(25, 19).Value = 0.08470104432599844
$ws.Cells.Item(25, 20).Value = 0.06259678260450952
# Row 26
$ws.Cells.Item(26, 7).Value = 70.46836733333333
$ws.Cells.Item(26, 8).Value = 211.405102
$ws.Cells.Item(26, 9).Value = 0.3576142671846927
$ws.Cells.Item(26, 10).Value = 0.36513419197984
$ws.Cells.Item(26, 13).Value = 13.62041
$ws.Cells.Item(26, 14).Value = 40.86123000000001
$ws.Cells.Item(26, 15).Value = 0.1369489801107134
$ws.Cells.Item(26, 16).Value = 0.1486879047638899
$ws.Cells.Item(26, 17).Value = 959.8080551106068
$ws.Cells.Item(26, 18).Value = 8638.272495995461
$ws.Cells.Item(26, 19).Value = 0.05897490916398383
$ws.Cells.Item(26, 20).Value = 0.05429103796313834

Write-Output "Done updating cells."